# "Add files via upload" - text of the color-legend cells in Sheet1 is
# updated to show each color's single-letter plotting code as a prefix,
# e.g. "Black" -> "k: Black", "Blue" -> "b: Blue", etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "k: Black"
$ws.Range("C11").Value = "b: Blue"
$ws.Range("C12").Value = "g: Green"
$ws.Range("C13").Value = "r: Red"
$ws.Range("C14").Value = "c: Cyan"
$ws.Range("C15").Value = "m: Magenta"
$ws.Range("C16").Value = "y: Yellow"
$ws.Range("C17").Value = "#XXXXXX: Arbitraty Color (#D185FF)"

# The last cell the author touched/selected in the sheet moved to C18
# (right below the last edited legend row).
$ws.Range("C18").Select()
